$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 124, shifting existing rows 124-162 down to 125-163
$ws.Rows.Item(124).EntireRow.Insert()

# Populate the newly inserted row 124 with the new weekly price record
$ws.Range("A124").Value = 7
$ws.Range("B124").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C124").Value = "Ñuble"
$ws.Range("D124").Value = 44463
$ws.Range("E124").Value = 16
$ws.Range("F124").Value = 100112009
$ws.Range("G124").Value = "Acelga"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 120
$ws.Range("K124").Value = 400
$ws.Range("L124").Value = 450
$ws.Range("M124").Value = 425
$ws.Range("N124").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O124").Value = "Provincia de Diguillín"
$ws.Range("P124").Value = 425
$ws.Range("Q124").Value = 1
$ws.Range("R124").Value = "Hortaliza"
